$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 24: label stays conceptually the same text ("...Per Iteration... (In Seconds)")
# but is now pointed at by the row that used to hold a different (batched) label.
$ws.Range("A24").Value = "Average Training Time Per Data Point Per Iteration For Neural Networks (In Seconds)"

# Row 25: speed label gains the "(On CPU)" qualifier.
$ws.Range("A25").Value = "Average Training Speed Per Data Point Per Iteration For Neural Networks (In Data Point Per Seconds) (On CPU)"

# Row 26: swaps to the "batched" time label, now qualified with "(On CPU)".
$ws.Range("A26").Value = "Average Training Time Per Data Point When Batched Per Iteration For Neural Networks (In Seconds) (On CPU)"

# Row 27: batched speed label, qualified with "(On CPU)".
$ws.Range("A27").Value = "Average Training Speed Per Data Point When Batched Per Iteration For Neural Networks (In Data Point Per Seconds) (On CPU)"

# Row 28: speedup factor label, qualified with "(On CPU)".
$ws.Range("A28").Value = "Speedup Factor Per Data Point When Batched For Neural Networks (On CPU)"

# Row 29: relative speedup factor label, qualified with "(On CPU)".
$ws.Range("A29").Value = "Relative Speedup Factor Per Data Point When Batched For Neural Networks When Compared With DataPredict Library (In Percentage) (On CPU)"

# Rows 25 and 26 now wrap onto a third line given the longer text, so their
# cached row heights grow from 30 to 45 (matching rows 27/29 which already
# wrap to 45).
$ws.Rows.Item(25).RowHeight = 45
$ws.Rows.Item(26).RowHeight = 45

# Move the sheet's active selection/cursor to A30.
[void]$ws.Range("A30").Select()
